$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write new strings in the order they first appear in the final workbook so the
# shared-strings table ends up in the same sequence as the target file:
#   interval, lb, ub, ws_check, ws_repair, switch_check
$ws.Range("B2").Value = "interval"
$ws.Range("C1").Value = "lb"
$ws.Range("D1").Value = "ub"
$ws.Range("B3").Value = "interval"
$ws.Range("B4").Value = "interval"

# Row 2: ws_fail
$ws.Range("A2").Value = "ws_fail"
$ws.Range("C2").Value = 0.0005
$ws.Range("D2").Value = 0.05
$ws.Range("E2").Value = $false

# Row 3: switch_fail
$ws.Range("A3").Value = "switch_fail"
$ws.Range("C3").Value = 0.00003
$ws.Range("D3").Value = 0.00125
$ws.Range("E3").Value = $false

# Row 4: line_fail
$ws.Range("A4").Value = "line_fail"
$ws.Range("C4").Value = 0.00002
$ws.Range("D4").Value = 0.001
$ws.Range("E4").Value = $false

# Row 5: ws_check (new)
$ws.Range("A5").Value = "ws_check"
$ws.Range("B5").Value = "interval"
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 3
$ws.Range("E5").Value = $false

# Row 6: ws_repair (new)
$ws.Range("A6").Value = "ws_repair"
$ws.Range("B6").Value = "interval"
$ws.Range("C6").Value = 0.5
$ws.Range("D6").Value = 1.2
$ws.Range("E6").Value = $false

# Row 7: switch_check (new)
$ws.Range("A7").Value = "switch_check"
$ws.Range("B7").Value = "interval"
$ws.Range("C7").Value = 8
$ws.Range("D7").Value = 12
$ws.Range("E7").Value = $false

# Column A width (mirrors the "AutoFit Column Width" the author applied; closest
# width this host's pixel-quantized ColumnWidth can reach to the recorded 11.8867)
$ws.Columns.Item(1).ColumnWidth = 11

# Selection
$ws.Range("E3").Select() | Out-Null
